$d = $word.ActiveDocument

function FindExecute($range, $findText, $replaceText) {
    return $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

function DeleteBetween($startText, $endText) {
    $r1 = $d.Content.Duplicate
    $r1.Find.Execute($startText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r2 = $d.Content.Duplicate
    $r2.Find.Execute($endText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $delRange = $d.Range($r1.Start, $r2.Start)
    $delRange.Text = ""
}

# ---------------------------------------------------------------------------
# Paragraph: first algorithm (compares a, b, c with ">" and writes a single
# winner immediately after the relevant "ise" branch instead of jumping to a
# separate "Yaz" step at the end).
# ---------------------------------------------------------------------------

# The find spans include the "x>y ." clause (which is wrapped in
# proofErr gramStart/gramEnd markers in the source) so that the replace
# collapses those now-meaningless grammar-check markers away, matching
# what Word does when it touches that text.
FindExecute $d.Content "a>b . a>c ise git 8 " "a>b . a>c ise yaz a ve git 8" | Out-Null
FindExecute $d.Content "b>a . b>c ise git 9" "b>a . b>c ise yaz b ve git 8" | Out-Null
FindExecute $d.Content "c>a . c>b ise git 10" "c>a . c>b ise yaz c ve git 8" | Out-Null

# Remove the now-obsolete separate "8-/9-/10- Yaz ... ve git 11" steps.
DeleteBetween "8-Yaz a ve git 11" "11-Dur"

FindExecute $d.Content "11-Dur" "8-Dur" | Out-Null

# ---------------------------------------------------------------------------
# Paragraph: second algorithm (orders a, b, c and writes the resulting order
# immediately after the relevant "ise" branch instead of jumping to a
# separate "Yaz" step at the end).
# ---------------------------------------------------------------------------

FindExecute $d.Content "a<b, b<c ise git 11" "a<b, b<c ise yaz a<b<c ve git 11" | Out-Null
FindExecute $d.Content "b<a, a<c ise git 12" "b<a, a<c ise yaz b<a<c ve git 11" | Out-Null
FindExecute $d.Content "a<c, c<b ise git 13" "a<c, c<b ise yaz a<c<b ve git 11" | Out-Null
FindExecute $d.Content "b<c, c<a ise git 14" "b<c, c<a ise yaz b<c<a ve git 11" | Out-Null
FindExecute $d.Content "c<a, a<b ise git 15" "c<a, a<b ise yaz c<a<b ve git 11" | Out-Null
FindExecute $d.Content "c<b, b<a ise git 16" "c<b, b<a ise yaz c<b<a ve git 11" | Out-Null

# Remove the now-obsolete separate "11-/.../16- Yaz ... git 17" steps.
DeleteBetween "11-Yaz a<b<c git 17" "17-Dur"

FindExecute $d.Content "17-Dur" "11-Dur" | Out-Null
